# Feria Lagunitas de Puerto Montt - Kiwi: weekly fruit/hortaliza price update.
# Inserts 3 new price rows (dated 2022-10-21 / serial 44855) right before the
# existing row 324, pushing the rest of the table down by 3 rows
# (old 324..401 -> new 327..404), then fills in the new rows' data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 324 (shifts rows 324:401 down to 327:404)
$ws.Rows(324).Resize(3).Insert()

# --- New row 324: Especial ---
$ws.Cells.Item(324, 1).Value = 4
$ws.Cells.Item(324, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(324, 3).Value = "Los Lagos"
$ws.Cells.Item(324, 4).Value = 44855
$ws.Cells.Item(324, 5).Value = 10
$ws.Cells.Item(324, 6).Value = "Fruta"
$ws.Cells.Item(324, 7).Value = 100101
$ws.Cells.Item(324, 8).Value = "Berries"
$ws.Cells.Item(324, 9).Value = 100101007
$ws.Cells.Item(324, 10).Value = "Kiwi"
$ws.Cells.Item(324, 11).Value = "Hayward"
$ws.Cells.Item(324, 12).Value = "Especial"
$ws.Cells.Item(324, 13).Value = 200
$ws.Cells.Item(324, 14).Value = 17000
$ws.Cells.Item(324, 15).Value = 17000
$ws.Cells.Item(324, 16).Value = 17000
$ws.Cells.Item(324, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(324, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(324, 19).Value = 1133
$ws.Cells.Item(324, 20).Value = 15

# --- New row 325: Primera ---
$ws.Cells.Item(325, 1).Value = 4
$ws.Cells.Item(325, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(325, 3).Value = "Los Lagos"
$ws.Cells.Item(325, 4).Value = 44855
$ws.Cells.Item(325, 5).Value = 10
$ws.Cells.Item(325, 6).Value = "Fruta"
$ws.Cells.Item(325, 7).Value = 100101
$ws.Cells.Item(325, 8).Value = "Berries"
$ws.Cells.Item(325, 9).Value = 100101007
$ws.Cells.Item(325, 10).Value = "Kiwi"
$ws.Cells.Item(325, 11).Value = "Hayward"
$ws.Cells.Item(325, 12).Value = "Primera"
$ws.Cells.Item(325, 13).Value = 200
$ws.Cells.Item(325, 14).Value = 15000
$ws.Cells.Item(325, 15).Value = 15000
$ws.Cells.Item(325, 16).Value = 15000
$ws.Cells.Item(325, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(325, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(325, 19).Value = 1000
$ws.Cells.Item(325, 20).Value = 15

# --- New row 326: Segunda ---
$ws.Cells.Item(326, 1).Value = 4
$ws.Cells.Item(326, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(326, 3).Value = "Los Lagos"
$ws.Cells.Item(326, 4).Value = 44855
$ws.Cells.Item(326, 5).Value = 10
$ws.Cells.Item(326, 6).Value = "Fruta"
$ws.Cells.Item(326, 7).Value = 100101
$ws.Cells.Item(326, 8).Value = "Berries"
$ws.Cells.Item(326, 9).Value = 100101007
$ws.Cells.Item(326, 10).Value = "Kiwi"
$ws.Cells.Item(326, 11).Value = "Hayward"
$ws.Cells.Item(326, 12).Value = "Segunda"
$ws.Cells.Item(326, 13).Value = 200
$ws.Cells.Item(326, 14).Value = 13000
$ws.Cells.Item(326, 15).Value = 13000
$ws.Cells.Item(326, 16).Value = 13000
$ws.Cells.Item(326, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(326, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(326, 19).Value = 867
$ws.Cells.Item(326, 20).Value = 15
